$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (class 1)
$ws.Range("B2").Value = 0.7955555555555556
$ws.Range("C2").Value = 0.9290657439446367
$ws.Range("D2").Value = 0.8571428571428572
$ws.Range("E2").Value = 1156

# Row 3 (class 2)
$ws.Range("B3").Value = 0.9525316455696202
$ws.Range("C3").Value = 0.9275808936825886
$ws.Range("D3").Value = 0.9398907103825137
$ws.Range("E3").Value = 649

# Row 4 (class 3)
$ws.Range("B4").Value = 0.8283950617283951
$ws.Range("C4").Value = 0.8515228426395939
$ws.Range("D4").Value = 0.8397997496871088
$ws.Range("E4").Value = 788

# Row 5 (class 4)
$ws.Range("B5").Value = 0.9662162162162162
$ws.Range("C5").Value = 0.4121037463976945
$ws.Range("D5").Value = 0.5777777777777778
$ws.Range("E5").Value = 347

# Row 6 (accuracy)
$ws.Range("B6").Value = 0.8469387755102041
$ws.Range("C6").Value = 0.8469387755102041
$ws.Range("D6").Value = 0.8469387755102041
$ws.Range("E6").Value = 0.8469387755102041

# Row 7 (macro avg)
$ws.Range("B7").Value = 0.8856746197674468
$ws.Range("C7").Value = 0.7800683066661284
$ws.Range("D7").Value = 0.8036527737475644

# Row 8 (weighted avg)
$ws.Range("B8").Value = 0.8591522434918055
$ws.Range("C8").Value = 0.8469387755102041
$ws.Range("D8").Value = 0.8377881991624915
